$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "Monday, Jan 09"
$ws.Range("C116").Value = "1:30 PM"
$ws.Range("D116").Value = "UNKNOWN"
$ws.Range("E116").Value = "Palanga"
$ws.Range("F116").Value = "(PLQ)"
$ws.Range("G116").Value = "Ryanair "
$ws.Range("H116").Value = "B738"
$ws.Range("I116").Value = "(SP-RSL)"
$ws.Range("J116").Value = "1:57 PM"
$ws.Range("L116").Value = "0 hours, 27 minutes"

$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "Monday, Jan 09"
$ws.Range("C117").Value = "1:40 PM"
$ws.Range("D117").Value = "FR3263"
$ws.Range("E117").Value = "Dublin"
$ws.Range("F117").Value = "(DUB)"
$ws.Range("G117").Value = "Ryanair "
$ws.Range("H117").Value = "B738"
$ws.Range("I117").Value = "(SP-RKI)"
$ws.Range("J117").Value = "1:54 PM"
$ws.Range("L117").Value = "0 hours, 14 minutes"

$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "Monday, Jan 09"
$ws.Range("C118").Value = "1:45 PM"
$ws.Range("D118").Value = "W61783"
$ws.Range("E118").Value = "Oslo"
$ws.Range("F118").Value = "(OSL)"
$ws.Range("G118").Value = "Wizz Air "
$ws.Range("H118").Value = "A320"
$ws.Range("I118").Value = "(HA-LWV)"
$ws.Range("J118").Value = "2:01 PM"
$ws.Range("L118").Value = "0 hours, 16 minutes"

$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "Monday, Jan 09"
$ws.Range("C119").Value = "2:00 PM"
$ws.Range("D119").Value = "FR3299"
$ws.Range("E119").Value = "Birmingham"
$ws.Range("F119").Value = "(BHX)"
$ws.Range("G119").Value = "Ryanair "
$ws.Range("H119").Value = "B738"
$ws.Range("I119").Value = "(SP-RKM)"
$ws.Range("J119").Value = "2:06 PM"
$ws.Range("L119").Value = "0 hours, 6 minutes"

$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "Monday, Jan 09"
$ws.Range("C120").Value = "2:15 PM"
$ws.Range("D120").Value = "LH1377"
$ws.Range("E120").Value = "Frankfurt"
$ws.Range("F120").Value = "(FRA)"
$ws.Range("G120").Value = "Lufthansa "
$ws.Range("H120").Value = "CRJ9"
$ws.Range("I120").Value = "(D-ACNF)"
$ws.Range("J120").Value = "3:18 PM"
$ws.Range("L120").Value = "1 hours, 3 minutes"

$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "Monday, Jan 09"
$ws.Range("C121").Value = "2:50 PM"
$ws.Range("D121").Value = "LO3838"
$ws.Range("E121").Value = "Warsaw"
$ws.Range("F121").Value = "(WAW)"
$ws.Range("G121").Value = "LOT "
$ws.Range("H121").Value = "E170"
$ws.Range("I121").Value = "(SP-LDI)"
$ws.Range("J121").Value = "2:50 PM"
$ws.Range("L121").Value = "0 hours, 0 minutes"

$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "Monday, Jan 09"
$ws.Range("C122").Value = "3:20 PM"
$ws.Range("D122").Value = "SK760"
$ws.Range("E122").Value = "Copenhagen"
$ws.Range("F122").Value = "(CPH)"
$ws.Range("G122").Value = "SAS "
$ws.Range("H122").Value = "A20N"
$ws.Range("I122").Value = "(EI-SIA)"
$ws.Range("J122").Value = "3:33 PM"
$ws.Range("L122").Value = "0 hours, 13 minutes"

$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "Monday, Jan 09"
$ws.Range("C123").Value = "3:20 PM"
$ws.Range("D123").Value = "W61733"
$ws.Range("E123").Value = "Stockholm"
$ws.Range("F123").Value = "(NYO)"
$ws.Range("G123").Value = "Wizz Air "
$ws.Range("H123").Value = "A320"
$ws.Range("I123").Value = "(HA-LYM)"
$ws.Range("J123").Value = "3:28 PM"
$ws.Range("L123").Value = "0 hours, 8 minutes"

$ws.Range("A124").Value = 123
$ws.Range("B124").Value = "Monday, Jan 09"
$ws.Range("C124").Value = "3:45 PM"
$ws.Range("D124").Value = "FR6127"
$ws.Range("E124").Value = "London"
$ws.Range("F124").Value = "(STN)"
$ws.Range("G124").Value = "Ryanair "
$ws.Range("H124").Value = "B738"
$ws.Range("I124").Value = "(SP-RKQ)"
$ws.Range("J124").Value = "3:56 PM"
$ws.Range("L124").Value = "0 hours, 11 minutes"

$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "Monday, Jan 09"
$ws.Range("C125").Value = "3:50 PM"
$ws.Range("D125").Value = "FR6099"
$ws.Range("E125").Value = "Stockholm"
$ws.Range("F125").Value = "(ARN)"
$ws.Range("G125").Value = "Ryanair "
$ws.Range("H125").Value = "B738"
$ws.Range("I125").Value = "(SP-RSW)"
$ws.Range("J125").Value = "4:04 PM"
$ws.Range("L125").Value = "0 hours, 14 minutes"

$ws.Range("A126").Value = 125
$ws.Range("B126").Value = "Monday, Jan 09"
$ws.Range("C126").Value = "4:35 PM"
$ws.Range("D126").Value = "W61791"
$ws.Range("E126").Value = "Larnaca"
$ws.Range("F126").Value = "(LCA)"
$ws.Range("G126").Value = "Wizz Air "
$ws.Range("H126").Value = "A321"
$ws.Range("I126").Value = "(HA-LTB)"
$ws.Range("J126").Value = "4:49 PM"
$ws.Range("L126").Value = "0 hours, 14 minutes"

$ws.Range("A127").Value = 126
$ws.Range("B127").Value = "Monday, Jan 09"
$ws.Range("C127").Value = "5:15 PM"
$ws.Range("D127").Value = "LO3816"
$ws.Range("E127").Value = "Warsaw"
$ws.Range("F127").Value = "(WAW)"
$ws.Range("G127").Value = "LOT "
$ws.Range("H127").Value = "E75S"
$ws.Range("I127").Value = "(SP-LIB)"
$ws.Range("J127").Value = "5:18 PM"
$ws.Range("L127").Value = "0 hours, 3 minutes"

